$d = $word.ActiveDocument

# 1) Add a first-line indent of 360 twips (18 pt) to the first paragraph.
$p1 = $d.Paragraphs(1)
$p1.FirstLineIndent = 18

# 2) Remove the _GoBack bookmark from the last paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
